# Auto-generated-assisted script to apply Goblin_Profits market data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1633.1632
$ws.Range("I15").Value = 1633.1632
$ws.Range("K15").Value = 4899.4896
$ws.Range("M15").Value = -4730.4896
$ws.Range("H74").Value = 4366.3335
$ws.Range("I74").Value = 3637.125
$ws.Range("J74").Value = 5824.75
$ws.Range("K74").Value = 3637.125
$ws.Range("L74").Value = 5824.75
$ws.Range("M74").Value = -2701.125
$ws.Range("N74").Value = -7696.75
$ws.Range("H75").Value = 375039420
$ws.Range("J75").Value = 375039420
$ws.Range("L75").Value = 375039420
$ws.Range("N75").Value = -375041292
$ws.Range("H77").Value = 4366.3335
$ws.Range("I77").Value = 3637.125
$ws.Range("J77").Value = 5824.75
$ws.Range("K77").Value = 18185.625
$ws.Range("L77").Value = 29123.75
$ws.Range("M77").Value = -13505.625
$ws.Range("N77").Value = -38483.75
$ws.Range("H78").Value = 375039420
$ws.Range("J78").Value = 375039420
$ws.Range("L78").Value = 1125118260
$ws.Range("N78").Value = -1125127620
$ws.Range("H106").Value = 2688.9333
$ws.Range("I106").Value = 2195.4167
$ws.Range("K106").Value = 2195.4167
$ws.Range("M106").Value = -1564.4167
$ws.Range("H110").Value = 39743.875
$ws.Range("J110").Value = 39743.875
$ws.Range("L110").Value = 39743.875
$ws.Range("N110").Value = -47923.875
$ws.Range("H116").Value = 15399.8
$ws.Range("J116").Value = 15999.667
$ws.Range("L116").Value = 15999.667
$ws.Range("N116").Value = -22883.667
$ws.Range("H132").Value = 1687.7273
$ws.Range("I132").Value = 1523.2333
$ws.Range("J132").Value = 3332.6667
$ws.Range("K132").Value = 4569.699900000001
$ws.Range("L132").Value = 9998.000100000001
$ws.Range("M132").Value = -2039.699900000001
$ws.Range("N132").Value = -15058.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4765.6567
$ws.Range("I32").Value = 2857.131
$ws.Range("K32").Value = 2857.131
$ws.Range("M32").Value = -2570.131
$ws.Range("H45").Value = 6145.409
$ws.Range("I45").Value = 6109.95
$ws.Range("K45").Value = 6109.95
$ws.Range("M45").Value = -5732.95
$ws.Range("H61").Value = 7418.35
$ws.Range("I61").Value = 6962.9414
$ws.Range("K61").Value = 6962.9414
$ws.Range("M61").Value = -6750.9414
$ws.Range("H74").Value = 2889.7222
$ws.Range("I74").Value = 2807.7856
$ws.Range("K74").Value = 2807.7856
$ws.Range("M74").Value = -1933.7856
$ws.Range("H77").Value = 2889.7222
$ws.Range("I77").Value = 2807.7856
$ws.Range("K77").Value = 14038.928
$ws.Range("M77").Value = -9670.928
$ws.Range("H110").Value = 4792.154
$ws.Range("I110").Value = 4792.154
$ws.Range("K110").Value = 4792.154
$ws.Range("M110").Value = -2747.154
$ws.Range("H122").Value = 5297689.5
$ws.Range("I122").Value = 7938152.5
$ws.Range("K122").Value = 23814457.5
$ws.Range("M122").Value = -23812007.5
$ws.Range("H136").Value = 7418.35
$ws.Range("I136").Value = 6962.9414
$ws.Range("K136").Value = 20888.8242
$ws.Range("M136").Value = -18338.8242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1649.258
$ws.Range("I94").Value = 1477.1818
$ws.Range("K94").Value = 1477.1818
$ws.Range("M94").Value = -1026.1818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 446.8889
$ws.Range("I22").Value = 377.875
$ws.Range("K22").Value = 377.875
$ws.Range("M22").Value = -27.875
$ws.Range("H31").Value = 4893.606
$ws.Range("I31").Value = 2442.3845
$ws.Range("J31").Value = 6486.9
$ws.Range("K31").Value = 2442.3845
$ws.Range("L31").Value = 6486.9
$ws.Range("M31").Value = -2147.3845
$ws.Range("N31").Value = -7076.9
$ws.Range("H34").Value = 4893.606
$ws.Range("I34").Value = 2442.3845
$ws.Range("J34").Value = 6486.9
$ws.Range("K34").Value = 2442.3845
$ws.Range("L34").Value = 6486.9
$ws.Range("M34").Value = -2240.3845
$ws.Range("N34").Value = -6890.9
$ws.Range("H43").Value = 24899.25
$ws.Range("J43").Value = 25666.666
$ws.Range("L43").Value = 25666.666
$ws.Range("N43").Value = -26034.666
$ws.Range("H101").Value = 24899.25
$ws.Range("J101").Value = 25666.666
$ws.Range("L101").Value = 25666.666
$ws.Range("N101").Value = -32156.666
$ws.Range("H105").Value = 4521.3335
$ws.Range("I105").Value = 3582
$ws.Range("J105").Value = 6400
$ws.Range("K105").Value = 3582
$ws.Range("L105").Value = 6400
$ws.Range("M105").Value = -1835
$ws.Range("N105").Value = -9894
$ws.Range("H133").Value = 38497.5
$ws.Range("I133").Value = 15000
$ws.Range("J133").Value = 61995
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 61995
$ws.Range("M133").Value = -12470
$ws.Range("N133").Value = -67055
$ws.Range("H141").Value = 226333.17
$ws.Range("J141").Value = 226333.17
$ws.Range("L141").Value = 226333.17
$ws.Range("N141").Value = -236693.17

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 3356.3333
$ws.Range("I25").Value = 4810
$ws.Range("J25").Value = 449
$ws.Range("K25").Value = 14430
$ws.Range("L25").Value = 1347
$ws.Range("M25").Value = -14261
$ws.Range("N25").Value = -1685
$ws.Range("H30").Value = 3356.3333
$ws.Range("I30").Value = 4810
$ws.Range("J30").Value = 449
$ws.Range("K30").Value = 14430
$ws.Range("L30").Value = 1347
$ws.Range("M30").Value = -14328
$ws.Range("N30").Value = -1551

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 49138.715
$ws.Range("J45").Value = 49138.715
$ws.Range("L45").Value = 49138.715
$ws.Range("N45").Value = -50256.715
$ws.Range("H70").Value = 30308614
$ws.Range("I70").Value = 41669344
$ws.Range("K70").Value = 41669344
$ws.Range("M70").Value = -41669074
$ws.Range("H73").Value = 30308614
$ws.Range("I73").Value = 41669344
$ws.Range("K73").Value = 41669344
$ws.Range("M73").Value = -41668408
$ws.Range("H80").Value = 6408.909
$ws.Range("I80").Value = 3999.3333
$ws.Range("K80").Value = 3999.3333
$ws.Range("M80").Value = -3001.3333
$ws.Range("H83").Value = 6408.909
$ws.Range("I83").Value = 3999.3333
$ws.Range("K83").Value = 19996.6665
$ws.Range("M83").Value = -15004.6665
$ws.Range("H97").Value = 1909.6316
$ws.Range("I97").Value = 1781.3529
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1781.3529
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1285.3529
$ws.Range("N97").Value = -3992
$ws.Range("H102").Value = 2605
$ws.Range("I102").Value = 806.6667
$ws.Range("K102").Value = 806.6667
$ws.Range("M102").Value = 815.3333
$ws.Range("H126").Value = 2671.5
$ws.Range("I126").Value = 2695.6667
$ws.Range("J126").Value = 2599
$ws.Range("K126").Value = 8087.000100000001
$ws.Range("L126").Value = 7797
$ws.Range("M126").Value = -5617.000100000001
$ws.Range("N126").Value = -12737
$ws.Range("H136").Value = 19273.926
$ws.Range("J136").Value = 19273.926
$ws.Range("L136").Value = 57821.778
$ws.Range("N136").Value = -62921.778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3801.75
$ws.Range("I61").Value = 2123
$ws.Range("K61").Value = 2123
$ws.Range("M61").Value = -1921
$ws.Range("H93").Value = 4850.033
$ws.Range("I93").Value = 1958.25
$ws.Range("K93").Value = 1958.25
$ws.Range("M93").Value = -710.25
$ws.Range("H113").Value = 3801.75
$ws.Range("I113").Value = 2123
$ws.Range("K113").Value = 2123
$ws.Range("M113").Value = 47
$ws.Range("H136").Value = 24785.375
$ws.Range("I136").Value = 4875.1113
$ws.Range("K136").Value = 14625.3339
$ws.Range("M136").Value = -12075.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1304.5
$ws.Range("I100").Value = 592
$ws.Range("J100").Value = 1779.5
$ws.Range("K100").Value = 1184
$ws.Range("L100").Value = 3559
$ws.Range("M100").Value = -643
$ws.Range("N100").Value = -4641
$ws.Range("H136").Value = 5503.3237
$ws.Range("I136").Value = 4265.92
$ws.Range("K136").Value = 12797.76
$ws.Range("M136").Value = -10247.76
$ws.Range("H137").Value = 96666.664
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
